# Commit: "when mapping to cdes search cde variables by pathology name and get latest version"
#
# 1) Add a new workbook-scoped (sheet-local) defined name that duplicates the
#    existing _FilterDatabase range, named "_xlnm._FilterDatabase_0".
# 2) Rewrite every "/root/..." concept-path string in column J (conceptPath)
#    of Sheet1 to "/dementia/...".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1) New defined name, scoped to Sheet1, referring to the same range as
#        the existing (Auto)Filter database A1:M180.
$filterRange = $ws.Range("A1:M180")
$ws.Names.Add("_xlnm._FilterDatabase_0", $filterRange)

# --- 2) Replace the "/root/" prefix with "/dementia/" across the whole
#        conceptPath column (J), rows 2..180 (row 1 is the header).
$lastRow = $ws.UsedRange.Rows.Count
$conceptCol = 10  # column J

for ($i = 2; $i -le $lastRow; $i++) {
    $cell = $ws.Cells.Item($i, $conceptCol)
    $text = $cell.Text
    if ($text -like "/root/*") {
        $cell.Value = $text -replace "^/root/", "/dementia/"
    }
}
